$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H20").Value = 12510.5
$ws.Range("I20").Value = 12510.5
$ws.Range("K20").Value = 12510.5
$ws.Range("M20").Value = -12280.5
$ws.Range("H29").Value = 1862.1936
$ws.Range("I29").Value = 80
$ws.Range("J29").Value = 2126.2222
$ws.Range("K29").Value = 240
$ws.Range("L29").Value = 6378.6666
$ws.Range("M29").Value = 41
$ws.Range("N29").Value = -6940.6666
$ws.Range("H31").Value = 178.75
$ws.Range("I31").Value = 178.75
$ws.Range("K31").Value = 536.25
$ws.Range("M31").Value = -306.25
$ws.Range("H34").Value = 5708
$ws.Range("I34").Value = 5124.615
$ws.Range("J34").Value = 9500
$ws.Range("K34").Value = 5124.615
$ws.Range("L34").Value = 9500
$ws.Range("M34").Value = -4921.615
$ws.Range("N34").Value = -9906
$ws.Range("H35").Value = 12510.5
$ws.Range("I35").Value = 12510.5
$ws.Range("K35").Value = 12510.5
$ws.Range("M35").Value = -12131.5
$ws.Range("H36").Value = 5708
$ws.Range("I36").Value = 5124.615
$ws.Range("J36").Value = 9500
$ws.Range("K36").Value = 5124.615
$ws.Range("L36").Value = 9500
$ws.Range("M36").Value = -4409.615
$ws.Range("N36").Value = -10930
$ws.Range("H127").Value = 62503012
$ws.Range("I127").Value = 250000450
$ws.Range("J127").Value = 3863.6667
$ws.Range("K127").Value = 750001350
$ws.Range("L127").Value = 11591.0001
$ws.Range("M127").Value = -749996390
$ws.Range("N127").Value = -21511.0001
$ws.Range("H132").Value = 20577136
$ws.Range("I132").Value = 2179674.5
$ws.Range("J132").Value = 333333980
$ws.Range("K132").Value = 6539023.5
$ws.Range("L132").Value = 1000001940
$ws.Range("M132").Value = -6536493.5
$ws.Range("N132").Value = -1000007000

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1236810.2
$ws.Range("I132").Value = 2327.4348
$ws.Range("J132").Value = 5292968
$ws.Range("K132").Value = 6982.3044
$ws.Range("L132").Value = 15878904
$ws.Range("M132").Value = -4452.3044
$ws.Range("N132").Value = -15883964

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H25").Value = 260666.67
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 260666.67
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 782000.01
$ws.Range("M25").ClearContents()
$ws.Range("N25").Value = -782338.01
$ws.Range("H29").Value = 309.44446
$ws.Range("I29").Value = 279.25
$ws.Range("J29").Value = 333.6
$ws.Range("K29").Value = 837.75
$ws.Range("L29").Value = 1000.8
$ws.Range("M29").Value = -560.75
$ws.Range("N29").Value = -1554.8
$ws.Range("H30").Value = 260666.67
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 260666.67
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = 782000.01
$ws.Range("M30").ClearContents()
$ws.Range("N30").Value = -782204.01
$ws.Range("H31").Value = 3500
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 3500
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 10500
$ws.Range("M31").ClearContents()
$ws.Range("N31").Value = -11076
$ws.Range("H35").Value = 1829.4117
$ws.Range("J35").Value = 1829.4117
$ws.Range("L35").Value = 5488.2351
$ws.Range("N35").Value = -6064.2351
$ws.Range("H36").Value = 2066.6667
$ws.Range("I36").Value = 200
$ws.Range("K36").Value = 600
$ws.Range("M36").Value = -431
$ws.Range("H139").Value = 2099.44
$ws.Range("I139").Value = 1675.6471
$ws.Range("J139").Value = 3000
$ws.Range("K139").Value = 5026.9413
$ws.Range("L139").Value = 9000
$ws.Range("M139").Value = 113.0587000000005
$ws.Range("N139").Value = -19280

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H29").Value = 8000
$ws.Range("J29").Value = 8000
$ws.Range("L29").Value = 8000
$ws.Range("N29").Value = -8580
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("M31").ClearContents()
$ws.Range("N31").ClearContents()
$ws.Range("H35").Value = 2999
$ws.Range("I35").Value = 2999
$ws.Range("K35").Value = 2999
$ws.Range("M35").Value = -2701
$ws.Range("H36").Value = 2413.2222
$ws.Range("I36").Value = 1333.3334
$ws.Range("J36").Value = 2953.1667
$ws.Range("K36").Value = 1333.3334
$ws.Range("L36").Value = 2953.1667
$ws.Range("M36").Value = -848.3334
$ws.Range("N36").Value = -3923.1667
$ws.Range("H37").Value = 0
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("M37").ClearContents()
$ws.Range("N37").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 44005
$ws.Range("I4").Value = 8000
$ws.Range("K4").Value = 8000
$ws.Range("M4").Value = -7887
$ws.Range("H28").Value = 44005
$ws.Range("I28").Value = 8000
$ws.Range("K28").Value = 8000
$ws.Range("M28").Value = -7768
$ws.Range("H29").Value = 14200
$ws.Range("I29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("M29").ClearContents()
$ws.Range("H30").Value = 26200
$ws.Range("I30").Value = 1000
$ws.Range("J30").Value = 32500
$ws.Range("K30").Value = 1000
$ws.Range("L30").Value = 32500
$ws.Range("M30").Value = -892
$ws.Range("N30").Value = -32716
$ws.Range("H31").Value = 1943.4546
$ws.Range("I31").Value = 499.33334
$ws.Range("J31").Value = 2485
$ws.Range("K31").Value = 499.33334
$ws.Range("L31").Value = 2485
$ws.Range("M31").Value = -251.33334
$ws.Range("N31").Value = -2981
$ws.Range("H35").Value = 3933.3333
$ws.Range("I35").Value = 3933.3333
$ws.Range("K35").Value = 3933.3333
$ws.Range("M35").Value = -3597.3333
$ws.Range("H37").Value = 44005
$ws.Range("I37").Value = 8000
$ws.Range("K37").Value = 8000
$ws.Range("M37").Value = -7893
$ws.Range("H97").Value = 24666.666
$ws.Range("J97").Value = 24666.666
$ws.Range("L97").Value = 24666.666
$ws.Range("N97").Value = -26648.666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 8180
$ws.Range("I29").Value = 6966.6665
$ws.Range("J29").Value = 10000
$ws.Range("K29").Value = 6966.6665
$ws.Range("L29").Value = 10000
$ws.Range("M29").Value = -6676.6665
$ws.Range("N29").Value = -10580
$ws.Range("H30").Value = 9670
$ws.Range("J30").Value = 9670
$ws.Range("L30").Value = 9670
$ws.Range("N30").Value = -9884
$ws.Range("H33").Value = 17330
$ws.Range("J33").Value = 17330
$ws.Range("L33").Value = 17330
$ws.Range("N33").Value = -17830
$ws.Range("H36").Value = 17330
$ws.Range("J36").Value = 17330
$ws.Range("L36").Value = 17330
$ws.Range("N36").Value = -17830
$ws.Range("H37").Value = 53352.668
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 53352.668
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 53352.668
$ws.Range("M37").ClearContents()
$ws.Range("N37").Value = -53758.668
$ws.Range("H132").Value = 3712.4827
$ws.Range("I132").Value = 3681.0435
$ws.Range("J132").Value = 3833
$ws.Range("K132").Value = 11043.1305
$ws.Range("L132").Value = 11499
$ws.Range("M132").Value = -8513.130500000001
$ws.Range("N132").Value = -16559
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()
